$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Active Cases"

$data = @(
    @('3398 BlueCross Elly Kay Mordialloc', 28),
    @('3564 Waverley Valley Aged Care Glen Waverley', 12),
    @('3601 Baptcare Westhaven community', 13),
    @('3647 Aurrum Aged Care Reservoir', 16),
    @('3653 Fronditha Thalpori St Albans Aged Care', 20),
    @('3975 Aurrum Aged Care Brunswick West', 13),
    @('4257 BlueCross The Gables Camberwell', 16),
    @('4295 Hope Aged Care Sunshine West', 15),
    @('4314 Estia Health Ardeer', 18),
    @('44095 Myrniong Primary School Myrniong', 17),
    @('44404 Castlemaine North Primary School Castlemaine', 60),
    @('44490 Armadale Primary School Armadale', 26),
    @('44507 Inverloch Primary School Inverloch', 11),
    @('44593 Torquay P-6 College Torquay', 26),
    @('44620 Canterbury Primary School Canterbury', 11),
    @('44623 Brunswick North West Primary School Brunswick West', 31),
    @('44642 Irymple South Primary School Irymple South', 10),
    @('44745 Briar Hill Primary School Briar Hill', 23),
    @('44765 Strathmore Primary School Strathmore', 10),
    @('4479 Whittlesea Lodge Whittlesea', 12),
    @('44799 Eastwood Primary School Ringwood East', 36),
    @('44855 Beverley Hills Primary School Doncaster East', 10),
    @('44893 Greenhills Primary School Greensborough', 14),
    @('44960 Thomastown West Primary School', 18),
    @('45013 Gladstone Views Primary School', 38),
    @('45147 Maramba Primary School Narre Warren', 11),
    @('45168 Ranfurly Primary School Mildura', 35),
    @('45257 Roxburgh Rise Primary School Roxburgh Park', 11),
    @('45305 Lockington Consolidated School Lockington', 41),
    @('4574 Village Glen Aged Care Residences Mornington', 10),
    @('45757 Saint Joseph''s Primary School Warragul', 13),
    @('45764 Our Lady Help of Christian''s Primary School Brunswick East', 11),
    @('45858 St Bernard''s Primary Coburg', 23),
    @('45958 Ave Maria College Aberfeldie Workplace', 21),
    @('46074 St Justin''s Catholic Primary School Wheelers Hill', 14),
    @('46078 Corpus Christi Primary School Werribee', 26),
    @('46086 St Kevin''s Primary School Hampton Park', 14),
    @('46208 Mount Scopus Memorial College Gandel Campus Burwood', 11),
    @('46320 St Mary''s Coptic Orthodox College Coolaroo', 10),
    @('46327 Victory Christian College Strathdale', 11),
    @('50279 Dallas Brooks Community Primary School Dallas', 11),
    @('50412 Geelong Grammar School Bostock House Campus Newtown', 13),
    @('52390 Our Lady of the Way Catholic Primary School Wallan', 47),
    @('52694 Pakenham Primary School Pakenham', 21),
    @('52777 Mirripoa Primary School Mount Duneed School Camp', 14),
    @('Australian Radio Network Richmond', 13),
    @('Cardinia Waters Retirement Village Pakenham', 16),
    @('Confirmed Omicron Sircuit Bar Fitzroy', 15),
    @('Confirmed Omicron Variant The Peel Hotel Collingwood', 20),
    @('Goodstart Early Learning Preston', 10),
    @('Green Gables Lodge Warburton', 17),
    @('Greendale Hotel Greendale', 23),
    @('JBS Australia Brooklyn', 41),
    @('Kororoit Creek Primary School Burnside Heights Oct-Dec', 21),
    @('PGL Camp Rumbug Foster North', 48),
    @('Rosebud Primary School Rosebud', 10),
    @('Social Gathering 11 Dec Windsor', 23),
    @('Social Gathering 11 December Fitzroy', 10),
    @('St Brigid''s Parish Primary School Mordialloc', 14),
    @('St Mary''s Parish Primary School', 12),
    @('StarTrack- Melbourne Tullamarine', 24),
    @('The George Lounge St Kilda', 16),
    @('Thomastown West Primary School Camp Doxa''s Malmsbury', 17)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
